$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arrays")
$ws.Activate()

# ---- Row 8: LC 1758 / Measy / Maximum Changes to Binary ----
$ws.Range("A8").Value = "LC 1758"
$ws.Range("B8").Value = "Measy"
$ws.Range("C8").Value = "Maximum Changes to Binary"

$text118 = @"
 * - The function iterates over each character in the string 's'.
 * - It counts the number of changes needed to make the string follow two patterns: "010101..." and "101010...".
 * - 'count0' is used to track the number of changes needed for the "010101..." pattern.
 *   - For even indices (0, 2, 4, ...), the expected character is '0'.
 *   - For odd indices (1, 3, 5, ...), the expected character is '1'.
 * - 'count1' is used to track the number of changes needed for the "101010..." pattern.
 *   - For even indices, the expected character is '1'.
 *   - For odd indices, the expected character is '0'.
 * - The function returns the minimum of 'count0' and 'count1', indicating the least number of changes needed.
"@
$text118 = $text118.TrimEnd("`r", "`n")
$ws.Range("E8").WrapText = $true
$ws.Range("E8").Value = $text118

# ---- Row 9: LC  / Medium / Valid Soduku ----
$ws.Range("A9").Value = "LC "
$ws.Range("B9").Value = "Medium"
$ws.Range("C9").Value = "Valid Soduku"

$text121 = @"
// Approach:
// 1. Use three 2D arrays (rows, colm, and box) to track the numbers present in each row, column, and 3x3 sub-box.
// 2. Iterate through each cell in the board.
//    - Convert the character at board[i][j] to an integer index (0-8) if it's not '.'.
//    - Calculate the index of the 3x3 sub-box.
//    - Check if the number already exists in the corresponding row, column, and sub-box.
//      - If it does, return false (invalid Sudoku).
//    - Otherwise, mark the number as seen in the corresponding row, column, and sub-box.
// 3. If no duplicates are found, return true (valid Sudoku).
"@
$text121 = $text121.TrimEnd("`r", "`n")
$ws.Range("E9").WrapText = $true
$ws.Range("E9").Value = $text121

# Row heights as computed/autofitted by Excel for the wrapped, multi-line content
$ws.Rows.Item(8).RowHeight = 187
$ws.Rows.Item(9).RowHeight = 170

# View state: scroll so row 6 is at top and select the last-edited cell (E8)
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E8").Select()
